# Atualizado por script em 08-11-2023 20:45
#
# The source scrape re-ran and reshuffled a handful of duplicate-date
# fixtures (rows 68/69, 79/80, 82/83, 98/99 each swap places with their
# neighbour) and appended one new match row (111) that had been missing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($sheet, $r1, $r2) {
    $range1 = $sheet.Range("F$r1`:V$r1")
    $range2 = $sheet.Range("F$r2`:V$r2")
    $v1 = $range1.Value()
    $v2 = $range2.Value()
    $range1.Value = $v2
    $range2.Value = $v1
}

# Rows 68/69, 79/80, 82/83 and 98/99 swap their home/away/odds/url content
# (columns F:V) while the index/country/tournament/season/date columns
# A:E stay put.
Swap-Rows $ws 68 69
Swap-Rows $ws 79 80
Swap-Rows $ws 82 83
Swap-Rows $ws 98 99

# Append the new, previously-missing match as row 111 (index 110).
# First clone the formatting of the last existing row so the bold/boxed
# index style (A) and the date number format (E) carry over correctly.
$ws.Range("A110:V110").Copy()
$ws.Range("A111:V111").PasteSpecial(-4122)  # xlPasteFormats
$ws.Application.CutCopyMode = 0

$ws.Range("A111").Value = 110
$ws.Range("B111").Value = "turkey"
$ws.Range("C111").Value = "super-lig"
$ws.Range("D111").Value = "2023-2024"
$ws.Range("E111").Value = 45238.75
$ws.Range("F111").Value = "Samsunspor"
$ws.Range("G111").Value = 2
$ws.Range("H111").Value = "Istanbulspor AS"
$ws.Range("I111").Value = 1
$ws.Range("J111").Value = 1.75
$ws.Range("K111").Value = "03/09/2023 18:14"
$ws.Range("L111").Value = 1.75
$ws.Range("M111").Value = "03/09/2023 18:14"
$ws.Range("N111").Value = 4.03
$ws.Range("O111").Value = "03/09/2023 18:14"
$ws.Range("P111").Value = 4.03
$ws.Range("Q111").Value = "03/09/2023 18:14"
$ws.Range("R111").Value = 4.79
$ws.Range("S111").Value = "03/09/2023 18:14"
$ws.Range("T111").Value = 4.79
$ws.Range("U111").Value = "03/09/2023 18:14"
$ws.Range("V111").Value = "https://www.betexplorer.com/football/turkey/super-lig/samsunspor-istanbulspor-as/EBGcs3V7/"
